$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Section_A")
$ws.Range('B2').Value = 'HS201'
$ws.Range('F2').Value = 'EC302'
$ws.Range('B3').Value = 'HS261 (Elective)'
$ws.Range('C3').Value = 'MA262'
$ws.Range('D3').Value = 'EC302'
$ws.Range('E3').Value = 'EC301'
$ws.Range('F3').Value = 'EC301'
$ws.Range('B5').Value = 'MA262'
$ws.Range('E5').Value = 'CS251 (Elective)'
$ws.Range('F5').Value = 'HS201'
$ws.Range('B6').Value = 'CS251 (Tutorial)'
$ws.Range('C6').Value = 'Free'
$ws.Range('E6').Value = 'Free'
$ws.Range('F6').Value = 'MA262 (Tutorial)'
$ws.Range('B7').Value = 'MA261'
$ws.Range('C7').Value = 'CS263'
$ws.Range('F7').Value = 'HS261 (Elective)'
$ws.Range('B8').Value = 'MA261 (Tutorial)'
$ws.Range('C8').Value = 'CS263 (Tutorial)'
$ws.Range('D8').Value = 'Free'
$ws.Range('E8').Value = 'EC302 (Tutorial)'
$ws.Range('F8').Value = 'HS261 (Tutorial)'

$ws = $wb.Worksheets.Item("Section_B")
$ws.Range('B2').Value = 'CS263'
$ws.Range('C2').Value = 'MA261'
$ws.Range('E2').Value = 'EC301'
$ws.Range('F2').Value = 'EC301'
$ws.Range('B3').Value = 'HS261 (Elective)'
$ws.Range('C3').Value = 'CS263'
$ws.Range('D3').Value = 'MA262'
$ws.Range('E3').Value = 'EC302'
$ws.Range('F3').Value = 'MA262'
$ws.Range('B5').Value = 'HS201'
$ws.Range('C5').Value = 'MA262'
$ws.Range('D5').Value = 'HS201'
$ws.Range('E5').Value = 'CS251 (Elective)'
$ws.Range('F5').Value = 'EC302'
$ws.Range('B6').Value = 'CS251 (Tutorial)'
$ws.Range('C6').Value = 'MA261 (Tutorial)'
$ws.Range('D6').Value = 'Free'
$ws.Range('E6').Value = 'EC302 (Tutorial)'
$ws.Range('C7').Value = 'HS201'
$ws.Range('D7').Value = 'EC302'
$ws.Range('F7').Value = 'HS261 (Elective)'
$ws.Range('B8').Value = 'MA262 (Tutorial)'
$ws.Range('C8').Value = 'Free'
$ws.Range('D8').Value = 'Free'
$ws.Range('E8').Value = 'CS263 (Tutorial)'
$ws.Range('F8').Value = 'HS261 (Tutorial)'

$ws = $wb.Worksheets.Item("Elective_Coordination")
$ws.Range('C2').Value = 'Mon'
$ws.Range('D3').Value = '15:30-17:00'
$ws.Range('C4').Value = 'Fri'
$ws.Range('C11').Value = 'Thu'
$ws.Range('D11').Value = '13:00-14:30'
$ws.Range('C13').Value = 'Mon'
$ws.Range('D13').Value = '14:30-15:30'
